# Apply data updates to "北京-漫展信息.xlsx"
# - 展览 (sheet1): G2 becomes text "不可售"; several F-column attendee counts increase.
# - 演出 (sheet2): several F-column attendee counts increase.
# - 本地生活 (sheet3): several F-column attendee counts increase.
# - 全部类型 (sheet4): several F-column attendee counts increase (mirrors the above sheets).

$wb = $excel.ActiveWorkbook

function Set-NumericCell {
    param(
        [object]$Worksheet,
        [string]$CellRef,
        [double]$NewValue
    )
    $Worksheet.Range($CellRef).Value = $NewValue
}

# ---------- Sheet: 展览 ----------
$ws1 = $wb.Worksheets.Item("展览")

# G2: was a plain number (55), now becomes the text "不可售"
$ws1.Range("G2").Value = "不可售"

Set-NumericCell $ws1 "F6"  5701
Set-NumericCell $ws1 "F11" 1547
Set-NumericCell $ws1 "F12" 10
Set-NumericCell $ws1 "F13" 25
Set-NumericCell $ws1 "F14" 650
Set-NumericCell $ws1 "F15" 1543
Set-NumericCell $ws1 "F16" 1543
Set-NumericCell $ws1 "F17" 1429
Set-NumericCell $ws1 "F18" 332
Set-NumericCell $ws1 "F19" 37
Set-NumericCell $ws1 "F20" 568
Set-NumericCell $ws1 "F21" 4096
Set-NumericCell $ws1 "F22" 4096
Set-NumericCell $ws1 "F23" 664
Set-NumericCell $ws1 "F25" 781
Set-NumericCell $ws1 "F27" 2241
Set-NumericCell $ws1 "F29" 319
Set-NumericCell $ws1 "F32" 1188
Set-NumericCell $ws1 "F35" 1101
Set-NumericCell $ws1 "F36" 1112

# ---------- Sheet: 演出 ----------
$ws2 = $wb.Worksheets.Item("演出")

Set-NumericCell $ws2 "F18" 278
Set-NumericCell $ws2 "F19" 200
Set-NumericCell $ws2 "F20" 481

# ---------- Sheet: 本地生活 ----------
$ws3 = $wb.Worksheets.Item("本地生活")

Set-NumericCell $ws3 "F4" 563
Set-NumericCell $ws3 "F5" 99
Set-NumericCell $ws3 "F6" 197

# ---------- Sheet: 全部类型 ----------
$ws4 = $wb.Worksheets.Item("全部类型")

Set-NumericCell $ws4 "F8"  563
Set-NumericCell $ws4 "F9"  99
Set-NumericCell $ws4 "F10" 5701
Set-NumericCell $ws4 "F21" 1547
Set-NumericCell $ws4 "F23" 10
Set-NumericCell $ws4 "F24" 25
Set-NumericCell $ws4 "F25" 1543
Set-NumericCell $ws4 "F27" 1429
Set-NumericCell $ws4 "F28" 332
Set-NumericCell $ws4 "F29" 37
Set-NumericCell $ws4 "F30" 568
Set-NumericCell $ws4 "F32" 4096
Set-NumericCell $ws4 "F33" 4096
Set-NumericCell $ws4 "F34" 664
Set-NumericCell $ws4 "F36" 781
Set-NumericCell $ws4 "F38" 2241
Set-NumericCell $ws4 "F44" 278
Set-NumericCell $ws4 "F45" 200
Set-NumericCell $ws4 "F46" 481
Set-NumericCell $ws4 "F49" 1101
Set-NumericCell $ws4 "F50" 1112
